# Second batch of data updates
# Discount Rate.xlsx - swap the US Social-Cost-of-Carbon discount-rate
# source/notes for the Mexican SHCP "Tasa Social de Descuento" source,
# bump the DR value from 3% to 10%, and make the DR sheet the active tab.

$wb  = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsDR    = $wb.Worksheets.Item("DR")

# --- "About" sheet: replace the source block (rows 3-7) -------------------
$wsAbout.Range("B3").Value = "Secretaría de Hacienda y Crédito Público"
$wsAbout.Range("B4").Value = 2021
$wsAbout.Range("B5").Value = "Tasa Social de Descuento (TSD)"

# B6 used to carry a clickable hyperlink to the old US source; drop the
# hyperlink (keeping the cell's existing styling) and point the plain text
# at the new Mexican source instead.
$wsAbout.Hyperlinks.Delete()
$wsAbout.Range("B6").Value = "https://www.gob.mx/shcp/documentos/tasa-social-de-descuento-tsd"

$wsAbout.Range("B7").Value = "Published 2015"

# Row 16 note: swap the old "3% discount rate" justification for the new
# "10% Mexican government rate" justification.
$wsAbout.Range("A16").Value = "We choose to use the official Mexican governmet 10% discount rate here, made consistent with"

# --- "DR" sheet: bump the discount rate from 3% to 10% --------------------
$wsDR.Range("B2").Value = 0.1

# --- Selection / active-tab bookkeeping (matches the authored workbook) ---
[void]$wsAbout.Range("A17").Select()
[void]$wsDR.Range("B3").Select()
$wsDR.Activate()
